$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (rank 0)
$ws.Range("D2").Value = '30.823.89'
$ws.Range("E2").Value = '  -1.20%  '

# Row 3 (rank 1)
$ws.Range("D3").Value = '1.940.16'
$ws.Range("E3").Value = '  -0.89%  '

# Row 4 (rank 2)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.21%  '

# Row 5 (rank 3)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.71'
$ws.Range("E5").Value = '  -1.73%  '

# Row 6 (rank 4)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.15%  '

# Row 7 (rank 5)
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4893'
$ws.Range("E7").Value = '  -0.03%  '

# Row 8 (rank 6)
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2930'
$ws.Range("E8").Value = '  -1.60%  '

# Row 9 (rank 7)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06890'
$ws.Range("E9").Value = '  +0.91%  '

# Row 10 (rank 8)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.27'
$ws.Range("E10").Value = '  -1.18%  '

# Row 11 (rank 9)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '105.00'
$ws.Range("E11").Value = '  -1.17%  '

# Row 12 (rank 10)
$ws.Range("D12").Value = '1.935.10'
$ws.Range("E12").Value = '  -0.95%  '

# Row 13 (rank 11)
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07751'
$ws.Range("E13").Value = '  -0.09%  '

# Row 14 (rank 12)
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.366'
$ws.Range("E14").Value = '  -1.09%  '

# Row 15 (rank 13)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6963'
$ws.Range("E15").Value = '  -2.95%  '

# Row 16 (rank 14)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '274.46'
$ws.Range("E16").Value = '  -3.02%  '

# Row 17 (rank 15)
$ws.Range("D17").Value = '30.819.41'
$ws.Range("E17").Value = '  -0.78%  '

# Row 18 (rank 16)
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007698'
$ws.Range("E18").Value = '  -0.82%  '

# Row 19 (rank 17)
$ws.Range("B19").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C19").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D19").Value = '2.205.15'
$ws.Range("E19").Value = '  +0.95%  '

# Row 20 (rank 18)
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.02'
$ws.Range("E20").Value = '  -1.49%  '

# Row 21 (rank 19)
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.09%  '

# Row 22 (rank 20)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.542'
$ws.Range("E22").Value = '  +0.79%  '

# Row 23 (rank 21)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.21%  '

# Row 24 (rank 22)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.555'
$ws.Range("E24").Value = '  -0.81%  '

# Row 25 (rank 23)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.790'
$ws.Range("E25").Value = '  -1.00%  '

# Row 26 (rank 24)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.74'
$ws.Range("E26").Value = '  -1.82%  '

# Row 27 (rank 25)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.61'
$ws.Range("E27").Value = '  -3.69%  '

# Row 28 (rank 26)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.154'
$ws.Range("E28").Value = '  -3.00%  '

# Row 29 (rank 27)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1038'
$ws.Range("E29").Value = '  -2.55%  '

# Row 30 (rank 28)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.394'
$ws.Range("E30").Value = '  -3.38%  '

# Row 31 (rank 29)
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.573'
$ws.Range("E31").Value = '  -2.93%  '

# Row 32 (rank 30)
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.554'
$ws.Range("E32").Value = '  -2.56%  '

# Row 33 (rank 31)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.359'
$ws.Range("E33").Value = '  -2.34%  '

# Row 34 (rank 32)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04884'
$ws.Range("E34").Value = '  -2.85%  '

# Row 35 (rank 33)
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7514'
$ws.Range("E35").Value = '  -1.46%  '

# Row 36 (rank 34)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.154'
$ws.Range("E36").Value = '  -0.88%  '

# Row 37 (rank 35)
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9998'
$ws.Range("E37").Value = '  -0.11%  '

# Row 38 (rank 36)
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.729'
$ws.Range("E38").Value = '  -0.29%  '

# Row 39 (rank 37)
$ws.Range("E39").Value = '  -2.72%  '

# Row 40 (rank 38)
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.574'
$ws.Range("E40").Value = '  +2.61%  '

# Row 41 (rank 39)
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.658'
$ws.Range("E41").Value = '  -1.87%  '

# Row 42 (rank 40)
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '78.24'
$ws.Range("E42").Value = '  +8.28%  '

# Row 43 (rank 41)
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.081'
$ws.Range("E43").Value = '  -4.95%  '

# Row 44 (rank 42)
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9043'
$ws.Range("E44").Value = '  +2.27%  '

# Row 45 (rank 43)
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4395'
$ws.Range("E45").Value = '  -3.27%  '

# Row 46 (rank 44)
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '107.65'
$ws.Range("E46").Value = '  -2.06%  '

# Row 47 (rank 45)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9989'
$ws.Range("E47").Value = '  -0.26%  '

# Row 48 (rank 46)
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.680'
$ws.Range("E48").Value = '  -0.16%  '

# Row 49 (rank 47)
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '979.93'
$ws.Range("E49").Value = '  +1.10%  '

# Row 50 (rank 48)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1240'
$ws.Range("E50").Value = '  -2.26%  '

# Row 51 (rank 49)
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.91'
$ws.Range("E51").Value = '  -0.11%  '
